$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 207, shifting existing rows 207:210 down to 208:211.
$ws.Rows(207).Insert()

# Populate the newly inserted row 207 with the new weekly data point.
$ws.Cells.Item(207, 1).Value = 5
$ws.Cells.Item(207, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(207, 3).Value = "Maule"
$ws.Cells.Item(207, 4).Value = 44448
$ws.Cells.Item(207, 5).Value = 7
$ws.Cells.Item(207, 6).Value = 100112043
$ws.Cells.Item(207, 7).Value = "Pepino ensalada"
$ws.Cells.Item(207, 8).Value = "Sin especificar"
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 300
$ws.Cells.Item(207, 11).Value = 16000
$ws.Cells.Item(207, 12).Value = 16000
$ws.Cells.Item(207, 13).Value = 16000
$ws.Cells.Item(207, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(207, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(207, 16).Value = 267
$ws.Cells.Item(207, 17).Value = 60
$ws.Cells.Item(207, 18).Value = "Hortaliza"
